# "incremento de verificação das jogadas"
# Expand the "jogo1" grid from 3x3 (A1:C3) to 5x5 (A1:E5), fixing up a
# couple of values that moved, and filling the newly-added rows/columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jogo1")

# Target 5x5 matrix of values (row-major, rows 1..5, cols A..E)
$values = @(
    @(0, -1, 0, 0, 0),
    @(0, 0, 0, 0, 0),
    @(-1, 0, 0, -1, 0),
    @(0, -1, 0, -1, 0),
    @(0, 0, 0, 0, 0)
)

for ($r = 1; $r -le 5; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}
